# Apply updated forecast vector data for this AR(2) winter series.
#
# A new earliest observation (2007 -> 2008) is inserted at the top of
# the table, shifting every existing row down by one, and a new most
# recent observation (2024 -> 2025) is appended at the bottom. The
# simulated y_0_forecast / y_1_forecast values are refreshed for every
# row to reflect the rerun of the bugfixed evaluator.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 2-18 down into rows 3-19, carrying their formatting along
# (this grows the table from 18 data rows, A1:E18, to 19 data rows,
# A1:E19) so every row keeps looking the same once its values are
# overwritten below.
$ws.Range("A2:E18").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

# The paste above created the brand-new row 19 using a generic default
# look; re-apply the same formatting used by the rest of the date
# column so it matches rows 2-18 exactly.
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Now write the refreshed values for every row.
$data = @()
$data += ,@(2, 39400, 2007, 7.226520411029069, 2008, $null)
$data += ,@(3, 39765, 2008, 4.268860212333636, 2009, $null)
$data += ,@(4, 40130, 2009, -7.266312015249776, 2010, $null)
$data += ,@(5, 40494, 2010, 6.958243460951929, 2011, 12.21658306395068)
$data += ,@(6, 40862, 2011, 9.469137444079934, 2012, 8.079264579851909)
$data += ,@(7, 41228, 2012, 3.358206407534947, 2013, 4.701432377325987)
$data += ,@(8, 41592, 2013, 0.3081076735359067, 2014, 3.972902167062387)
$data += ,@(9, 41957, 2014, 3.901355411819707, 2015, 4.658857392675264)
$data += ,@(10, 42321, 2015, 5.331683351557981, 2016, 4.089819750351786)
$data += ,@(11, 42689, 2016, 3.254758369308375, 2017, 2.313009565865753)
$data += ,@(12, 43053, 2017, 5.246209615995667, 2018, 4.784022165496182)
$data += ,@(13, 43418, 2018, 4.86255966374296, 2019, 4.112897401876747)
$data += ,@(14, 43783, 2019, 2.764740011159428, 2020, 1.643374185611401)
$data += ,@(15, 44159, 2020, -7.260793671746435, 2021, 0.00562230452727519)
$data += ,@(16, 44525, 2021, 4.097586525396268, 2022, 3.9116372951149)
$data += ,@(17, 44890, 2022, 7.824284864703746, 2023, 2.586378346096296)
$data += ,@(18, 45254, 2023, -1.24502235313334, 2024, -1.561801765212567)
$data += ,@(19, 45618, 2024, -1.735114423676209, 2025, 2.409056355286521)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    if ($row[5] -ne $null) {
        $ws.Cells.Item($r, 5).Value = $row[5]
    }
}
